$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Value = "x"
$ws.Range("E14").Value = "x"
$ws.Range("G14").Value = "x"
$ws.Range("H14").Value = "x"
$ws.Range("I14").Value = "x"
$ws.Range("J14").Value = "x"
$ws.Range("K14").Value = "x"
$ws.Range("L14").Value = "x"
$ws.Range("M14").Value = "x"

$ws.Range("M17").Select()
